# ---------------------------------------------------------------------------
# Commit: "added number of records updated and failed + switched to reason
# codes"
#
#  * Reason labels in column R (and the shared-string table) are switched
#    from long descriptions ("Discontinued/ Obsolete", "Sales Related") to
#    short reason codes ("O2", "RJ") and two more codes are introduced
#    ("MS", "SO") for the two freshly appended records.
#  * Two new data rows are appended (rows 4 and 5), each carrying a SKU code
#    (col C), a sales-document number (col N) and a reason code (col R).
#  * The sheet's "best fit" column now sits under the Sales-document column
#    (N) instead of the previously empty Sales-doc-type column (I), because
#    that's where the newly typed long numbers live.
#  * The active selection moves to the last entered cell, C5.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- switch the existing rows from long text reasons to short codes -------
$ws.Range("R2").Value = "O2"
$ws.Range("R3").Value = "RJ"

# --- append the two new records --------------------------------------------
$ws.Range("C4").Value = 1601002
$ws.Range("N4").Value = 102343878
$ws.Range("R4").Value = "MS"

$ws.Range("C5").Value = 1600544
$ws.Range("N5").Value = 102343878
$ws.Range("R5").Value = "SO"

# --- column N now holds the widest content, so it picks up the custom,
#     best-fit width that used to live on column I -------------------------
$ws.Columns.Item(14).ColumnWidth = 8.92

# --- leave the selection on the last cell entered, matching the author ----
$ws.Range("C5").Select()
